$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 28)
$ws.Range("A28").Value = "Neuse River Brewing"
$ws.Range("B28").Value = "Raleigh"
$ws.Range("C28").Value = "Neuse River Burger"
$ws.Range("D28").Value = "Brassiere/Burgers"
$ws.Range("E28").Value = 35.80457
$ws.Range("F28").Value = -78.6325

# Update the selection / view to match the author's saved state
$ws.Range("A28:F28").Select()
$excel.ActiveWindow.ScrollRow = 4
